$d = $word.ActiveDocument

# The document contains three "<id>...</id>" fields that were previously split
# across three separate runs (opening tag / id value / closing tag), each with
# their own distinct run formatting. The edit collapses each of these into a
# single run containing the full "<id>value</id>" text (using the formatting
# of the original "<id>" run), and renames/shortens the id values themselves.
#
# Using Find & Replace with exact text on the whole three-run span merges the
# matched runs into a single run that keeps the formatting of the first
# (leading) run of the match, which reproduces the target XML exactly.

$d.Content.Find.Execute("<id>p004v_a3</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p004v_3</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p005r_a1</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p005r_1</id>", 2) | Out-Null
$d.Content.Find.Execute("<id>p005r_a2</id>", $true, $false, $false, $false, $false, $true, 1, $false, "<id>p005r_2</id>", 2) | Out-Null
